# Resources.xlsx update
# - adds two new columns ("Unit" and "Notes" detail) and renames the old
#   descriptive "Notes" column to a short "Analog to" label
# - the long "analog to ..." sentences are split into a short keyword
#   (column C) + unit (column D) + trimmed explanatory note (column E)
# - updates the active-cell selection and column widths to match the
#   reworked table layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=1;  A='Resources'; B=1;     C='Analog to';                 D='Unit';                E='Notes' },
    @{ Row=2;  A='R1';  B=0.2;  C='population';                D='million people';      E=" the amount of people in a country is only a small indicator of the country's prosperity. Some small countries are very wealthy and some are very poor. Overall, does indicate some sense of wealth for a country." },
    @{ Row=3;  A='R2';  B=0.75; C='metallic elements';          D='million tons';        E='essential for metallic alloy creation and electronic creation ' },
    @{ Row=4;  A='R3';  B=0.5;  C='timber';                     D='million tons';        E='used in all forms of construction, but not a particularly rare resource' },
    @{ Row=5;  A='R4';  B=1;    C='available land';             D='million acres';       E='valued at twice the weight as water because land limits how much housing/farm/factories can be created that bring large amounts of prosperity. ' },
    @{ Row=6;  A='R5';  B=1;    C='renewable energy';           D='million kW';          E="renewable energy valued at 1 - in direct correlation to renewable energy waste's weight being -1" },
    @{ Row=7;  A='R6';  B=1.5;  C='fossil fuel energy';         D='million kW';          E='fossil fuels create more energy than green sources, but their waste is higher to indicate penalty for using nonrenewables.' },
    @{ Row=8;  A='R7';  B=0.5;  C='water';                      D='billion gallons';     E='0.5 chosen as the baseline for which all other raw resources are weighted. Essential for life and is involved in other types of resource creation, but is not rare.' },
    @{ Row=9;  A='R8';  B=0.5;  C='animals';                    D='million animals';     E='used for farms and food. Not particularly rare and has only a few use cases' },
    @{ Row=10; A='R9';  B=0.5;  C='plants';                     D='million tons';        E='used for farms and food, also produces fresh oxygen. Not particularly rare as well.' },
    @{ Row=11; A='R18'; B=3;    C='metallic alloys';            D='million tons';        E='weighted at 2 to account for -1 alloy waste weight. Alloy + alloy waste = 2 in weight, compared to 1.5 in lost input resources' },
    @{ Row=12; A='R19'; B=15;   C='housing';                    D='million homes';       E='weighted at 15 to account for -2 housing waste weight. Input resources lost have combined weight of 12.25.' },
    @{ Row=13; A='R20'; B=5;    C='electronics';                D='million gadgets';     E='weighted at 5 since 2 electronics and 1 waste is created - these total to 9 in weight compared to 8.25 of lost input resources' },
    @{ Row=14; A='R21'; B=1;    C='farm';                       D='million acres' },
    @{ Row=15; A='R22'; B=1;    C='factory';                    D='thousand factories' },
    @{ Row=16; A="R1'";  B=-1;  C='population waste' },
    @{ Row=17; A="R5'";  B=-1;  C='renewable energy waste';     E="renewable energies' waste weighted at -1 so that there is no net loss in using renewable energies" },
    @{ Row=18; A="R6'";  B=-2;  C='nonrenewable energy waste';  E='nonrenewable energy waste is weighted higher than the weight of nonrenewable energy, to discourage fossil fuel use' },
    @{ Row=19; A="R18'"; B=-1;  C='metallic alloys waste' },
    @{ Row=20; A="R19'"; B=-2;  C='housing waste' },
    @{ Row=21; A="R20'"; B=-1;  C='electronics waste' },
    @{ Row=22; A="R21'"; B=-1;  C='farm waste' },
    @{ Row=23; A="R22'"; B=-1;  C='factory waste' }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    if ($r.ContainsKey('D')) {
        $ws.Cells.Item($row, 4).Value = $r.D
    }
    if ($r.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = $r.E
    }
}

# Column widths - the new layout widens the descriptive columns so the
# table reads like a reference sheet (bestfit-style sizing).
$ws.Columns.Item(1).ColumnWidth = 8.5
$ws.Columns.Item(2).ColumnWidth = 5.67
$ws.Columns.Item(3).ColumnWidth = 22.5
$ws.Columns.Item(4).ColumnWidth = 14.33
$ws.Columns.Item(5).ColumnWidth = 159.5

# Selection moved to E10 (matches the new "Notes" column) and the window
# is restored to a plain top-left, maximized layout.
$ws.Range("E10").Select()
